$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record was added to the dataset. It is inserted as a
# new row 180, which pushes the previously-existing rows 180-239 down to
# rows 181-240 (their contents are left untouched by the insert itself).
$ws.Rows.Item(180).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Cells.Item(180, 1).Value = 3
$ws.Cells.Item(180, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(180, 3).Value = "Coquimbo"
$ws.Cells.Item(180, 4).Value = 44524
$ws.Cells.Item(180, 5).Value = 5
$ws.Cells.Item(180, 6).Value = 100112043
$ws.Cells.Item(180, 7).Value = "Pepino ensalada"
$ws.Cells.Item(180, 8).Value = "Sin especificar"
$ws.Cells.Item(180, 9).Value = "Primera"
$ws.Cells.Item(180, 10).Value = 105
$ws.Cells.Item(180, 11).Value = 7500
$ws.Cells.Item(180, 12).Value = 8000
$ws.Cells.Item(180, 13).Value = 7762
$ws.Cells.Item(180, 14).Value = "`$/caja 70 unidades"
$ws.Cells.Item(180, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(180, 16).Value = 111
$ws.Cells.Item(180, 17).Value = 70
$ws.Cells.Item(180, 18).Value = "Hortaliza"
